$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation for Puerro was recorded and inserted
# as row 38 ("Fruta / hortaliza, semanal"), pushing the existing rows
# 38-81 down to 39-82.
$ws.Rows(38).Insert()

$ws.Range("A38").Value = 9
$ws.Range("B38").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C38").Value = "Metropolitana"
$ws.Range("D38").Value = 44539
$ws.Range("E38").Value = 13
$ws.Range("F38").Value = 100112005
$ws.Range("G38").Value = "Puerro"
$ws.Range("H38").Value = "Sin especificar"
$ws.Range("I38").Value = "Primera"
$ws.Range("J38").Value = 133
$ws.Range("K38").Value = 6000
$ws.Range("L38").Value = 7000
$ws.Range("M38").Value = 6504
$ws.Range("N38").Value = "`$/paquete 20 unidades"
$ws.Range("O38").Value = "Provincia de Chacabuco"
$ws.Range("P38").Value = 325
$ws.Range("Q38").Value = 20
$ws.Range("R38").Value = "Hortaliza"
